# Add other racial groups into the national poverty table.
# - Insert two new data rows (pop_id 4 and 6) right after the existing
#   "United States" (pop_id 1) total row, shifting the remaining rows down.
# - Append two new data rows (pop_id 400 and 451) at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 3 (pushes old rows 3:68 down to 5:70).
$ws.Rows("3:4").Insert()

# New rows inserted near the top of the table.
$topRows = @(
    @(4, 38228744, 0, 10321254, 0, 27907492, 0),
    @(6, 2481414, 0, 702127, 0, 1779287, 0)
)

$r = 3
foreach ($row in $topRows) {
    $ws.Range("A$r").Value = "United States"
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $ws.Range("H$r").Value = $row[6]
    $r = $r + 1
}

# New rows appended at the bottom of the table.
$bottomRows = @(
    @(400, 53139880, 0, 12915617, 0, 40224264, 0),
    @(451, 192733728, 0, 20750472, 0, 171983264, 0)
)

$r = 71
foreach ($row in $bottomRows) {
    $ws.Range("A$r").Value = "United States"
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $ws.Range("H$r").Value = $row[6]
    $r = $r + 1
}
